# Auto-generated edit script applying numeric corrections to H:N
# (currentAveragePrice / LevePrice / LeveProfit columns) across several
# worksheets, per the scheduled-runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H2").Value = 325
$ws.Range("I2").Value = 340
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 340
$ws.Range("L2").Value = 250
$ws.Range("M2").Value = -227
$ws.Range("N2").Value = -476

$ws.Range("H17").Value = 1303.3334
$ws.Range("I17").Value = 700
$ws.Range("J17").Value = 1475.7142
$ws.Range("K17").Value = 2100
$ws.Range("L17").Value = 4427.142599999999
$ws.Range("M17").Value = -1932
$ws.Range("N17").Value = -4763.142599999999

$ws.Range("H132").Value = 4048
$ws.Range("I132").Value = 1667.4
$ws.Range("J132").Value = 9999.5
$ws.Range("K132").Value = 5002.200000000001
$ws.Range("L132").Value = 29998.5
$ws.Range("M132").Value = -2472.200000000001
$ws.Range("N132").Value = -35058.5

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H45").Value = 31508.285
$ws.Range("I45").Value = 8111.6
$ws.Range("J45").Value = 90000
$ws.Range("K45").Value = 8111.6
$ws.Range("L45").Value = 90000
$ws.Range("M45").Value = -7734.6
$ws.Range("N45").Value = -90754

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0

$ws.Range("H122").Value = 5550
$ws.Range("I122").Value = 4241.75
$ws.Range("J122").Value = 6596.6
$ws.Range("K122").Value = 12725.25
$ws.Range("L122").Value = 19789.8
$ws.Range("M122").Value = -10275.25
$ws.Range("N122").Value = -24689.8

$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0

$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0

$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0

$ws.Range("H122").Value = 90000
$ws.Range("I122").Value = 90000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 90000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -85100

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0

$ws.Range("H126").Value = 120000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 120000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 120000
$ws.Range("N126").Value = -129880

$ws.Range("H127").Value = 49750
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 49750
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 49750
$ws.Range("N127").Value = -59670

$ws.Range("H128").Value = 11664.333
$ws.Range("I128").Value = 11664.333
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 34992.999
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -32502.999

$ws.Range("H129").Value = 28916.666
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 28916.666
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 28916.666
$ws.Range("N129").Value = -38916.666

$ws.Range("H130").Value = 97173.75
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 97173.75
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 97173.75
$ws.Range("N130").Value = -107213.75

$ws.Range("H131").Value = 88266.664
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 88266.664
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 88266.664
$ws.Range("N131").Value = -98346.664

$ws.Range("H132").Value = 70000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 70000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 70000
$ws.Range("N132").Value = -80120

$ws.Range("H133").Value = 70000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 70000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120

$ws.Range("H134").Value = 3347.4167
$ws.Range("I134").Value = 3451.7273
$ws.Range("J134").Value = 2200
$ws.Range("K134").Value = 10355.1819
$ws.Range("L134").Value = 6600
$ws.Range("M134").Value = -7820.1819
$ws.Range("N134").Value = -11670

$ws.Range("H135").Value = 60780
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 60780
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 60780
$ws.Range("N135").Value = -70920

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0

$ws.Range("H138").Value = 15180
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 15180
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 15180
$ws.Range("N138").Value = -25460

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0

$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0

$ws.Range("H141").Value = 159977.5
$ws.Range("I141").Value = 159977.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 159977.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -154797.5

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H12").Value = 85.166664
$ws.Range("I12").Value = 55.5
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 166.5
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = 6.5
$ws.Range("N12").Value = -646

$ws.Range("H116").Value = 4568.421
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4568.421
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 13705.263
$ws.Range("N116").Value = -20589.263

$ws.Range("H132").Value = 10424
$ws.Range("I132").Value = 3508.8
$ws.Range("J132").Value = 45000
$ws.Range("K132").Value = 31579.2
$ws.Range("L132").Value = 405000
$ws.Range("M132").Value = -29049.2
$ws.Range("N132").Value = -410060

$ws.Range("H133").Value = 50
$ws.Range("I133").Value = 50
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 150
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = 4910

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0

$ws.Range("M134:N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H122").Value = 3428.375
$ws.Range("I122").Value = 3489.5715
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 10468.7145
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -8018.7145
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H46").Value = 4500
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 4375
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 4375
$ws.Range("M46").Value = -4812
$ws.Range("N46").Value = -4751

$ws.Range("H100").Value = 6330
$ws.Range("I100").Value = 7000
$ws.Range("J100").Value = 4990
$ws.Range("K100").Value = 7000
$ws.Range("L100").Value = 4990
$ws.Range("M100").Value = -6459
$ws.Range("N100").Value = -6072

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0

$ws.Range("H127").Value = 115000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 115000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 115000
$ws.Range("N127").Value = -124920

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0

$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0

$ws.Range("H130").Value = 90000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 90000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 90000
$ws.Range("N130").Value = -100040

$ws.Range("H131").Value = 78828
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 78828
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 78828
$ws.Range("N131").Value = -88908

$ws.Range("H132").Value = 3719.8
$ws.Range("I132").Value = 2616.6667
$ws.Range("J132").Value = 5374.5
$ws.Range("K132").Value = 7850.000100000001
$ws.Range("L132").Value = 16123.5
$ws.Range("M132").Value = -5320.000100000001
$ws.Range("N132").Value = -21183.5

$ws.Range("H133").Value = 130000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 130000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 130000
$ws.Range("N133").Value = -135060

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0

$ws.Range("H136").Value = 2824.5
$ws.Range("I136").Value = 2717.6365
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 8152.9095
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -5602.9095
$ws.Range("N136").Value = -17100

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0

$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H122").Value = 4050.8333
$ws.Range("I122").Value = 3749.75
$ws.Range("J122").Value = 4653
$ws.Range("K122").Value = 11249.25
$ws.Range("L122").Value = 13959
$ws.Range("M122").Value = -8799.25
$ws.Range("N122").Value = -18859

$ws.Range("H132").Value = 1610.25
$ws.Range("I132").Value = 1302.1818
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 3906.5454
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1376.5454
$ws.Range("N132").Value = -20057
